$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 3
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Select cell E3 (reflects <selection activeCell="E3" sqref="E3"/> in sheetView)
$ws.Range("E3").Select()
